$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New row 4 data: MSISDN / IMEI sample values ---------------------------
# C4 (IMEI-looking value) gets a distinct font (Segoe UI / FF242424) first so
# the new font-carrying cell style lands before the plain text style below.
$ws.Range("C4").NumberFormat = "@"
$ws.Range("C4").Font.Name = "Segoe UI"
$ws.Range("C4").Font.Color = 2368548
$ws.Range("C4").Value = "883333340718342"

# A4 / B4 (MSISDN-like values) use the plain/default font but as text.
$ws.Range("A4").NumberFormat = "@"
$ws.Range("A4").Value = "491453906"
$ws.Range("B4").NumberFormat = "@"
$ws.Range("B4").Value = "3043209863"

# --- Re-apply text format to the header / label rows ------------------------
$ws.Range("A1:C1").NumberFormat = "@"
$ws.Range("A3:C3").NumberFormat = "@"

# --- Row height for the new row ---------------------------------------------
$ws.Rows.Item(4).RowHeight = 16.5

# --- Column C a bit wider so the IMEI value fits -----------------------------
$ws.Columns.Item(3).ColumnWidth = 18

# --- Move the active selection -----------------------------------------------
[void]$ws.Range("D7").Select()

# --- Page setup: portrait orientation ----------------------------------------
$ws.PageSetup.Orientation = 1
